$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# The shared string used by column N's header was renamed from
# "Valores_consumo Valor" to "Valores_consumo". Re-assigning the cell's
# value rewrites the underlying shared string (and naturally reshuffles the
# shared-string table / indices exactly like Excel does).
$ws.Range("N1").Value = "Valores_consumo"

# The sheet view had scrolled so column H was left-most and M12 was the
# active/selected cell (previously B1 was left-most with I1 selected).
$win = $excel.ActiveWindow
$win.ScrollColumn = 8
$win.ScrollRow = 1
$ws.Range("M12").Select()
